$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45461 -> 45462, i.e. 2024-06-18 -> 2024-06-19) for every data row.
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45461) {
        $cell.Value2 = 45462
    }
}
